# Generate Report for Handback
# Updates the Overview / zh-cn / de-de sheets with handback results:
#  - Status moves from "In Translation" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File columns get populated
#  - Latest Handback DateTime gets a real timestamp (per locale)
#  - The "Latest Target File" column also becomes a hyperlink to the source .md

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a8ddfaeb05fc01e265c861548df4c29fbb8991/e2e/"

$file1 = "e2999384-cfac-441a-bc21-e91b789e10dd"
$file2 = "e905d38a-328b-4fab-b561-e615da70eee2"

$status = "Handed back: in sync with en-US"

# ---- Overview sheet: widen the per-locale status columns (E, F) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status

$wsZh.Range("I2").Value = "$file1.md"
$wsZh.Range("J2").Value = "$file1.c2898fcbec778b450488fc4b850b16b2fd1a01bb.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 22:24:39"

$wsZh.Range("I3").Value = "$file2.md"
$wsZh.Range("J3").Value = "$file2.e7a2d177ba6f10f69b3977a37c93b3d076633b6e.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-26 22:24:39"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $repoBase + "$file1.md", "", "", "$file1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $repoBase + "$file2.md", "", "", "$file2.md")

$wsZh.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status

$wsDe.Range("I2").Value = "$file1.md"
$wsDe.Range("J2").Value = "$file1.c2898fcbec778b450488fc4b850b16b2fd1a01bb.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 22:24:45"

$wsDe.Range("I3").Value = "$file2.md"
$wsDe.Range("J3").Value = "$file2.e7a2d177ba6f10f69b3977a37c93b3d076633b6e.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-26 22:24:45"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $repoBase + "$file1.md", "", "", "$file1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $repoBase + "$file2.md", "", "", "$file2.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.1666666666667
